# B6-PowerPoint.pptx — Wed, May 13, 2020  6:05:04 AM
#
# 1) Three tables (slides 14, 15, 16) switch from the deck's custom
#    "Table_0" style to the built-in "No Style, Table Grid" style.
# 2) The slide-master theme ("Integral" / Red Violet colour scheme) is
#    replaced with the stock default Office colour scheme (the colours
#    that used to live only on the notes-master's theme).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newStyleId = "{9BD3F43D-DBBE-4A4B-A4A6-231CD179C196}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newStyleId)
}

# --- 2. Swap the theme colour scheme over to the default Office one --
# (Office theme colours: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme

$colorScheme.Colors(1).RGB = 0          # dk1      000000
$colorScheme.Colors(2).RGB = 16777215   # lt1      FFFFFF
$colorScheme.Colors(3).RGB = 6968388    # dk2      44546A
$colorScheme.Colors(4).RGB = 15132391   # lt2      E7E6E6
$colorScheme.Colors(5).RGB = 13998939   # accent1  5B9BD5
$colorScheme.Colors(6).RGB = 3243501    # accent2  ED7D31
$colorScheme.Colors(7).RGB = 10855845   # accent3  A5A5A5
$colorScheme.Colors(8).RGB = 49407      # accent4  FFC000
$colorScheme.Colors(9).RGB = 12874308   # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477   # folHlink 954F72
